$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.861.57"
$ws.Range("E2").Value = "  +7.93%  "
$ws.Range("D3").Value = "'3.554.23"
$ws.Range("E3").Value = "  +10.48%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'192.76"
$ws.Range("E5").Value = "  +10.69%  "
$ws.Range("D6").Value = "'556.82"
$ws.Range("E6").Value = "  +8.41%  "
$ws.Range("D7").Value = "'3.547.38"
$ws.Range("E7").Value = "  +10.38%  "
$ws.Range("D8").Value = "'0.611"
$ws.Range("E8").Value = "  +3.56%  "
$ws.Range("E9").Value = "  -0.09%  "
$ws.Range("D10").Value = "'0.644"
$ws.Range("E10").Value = "  +8.20%  "
$ws.Range("D11").Value = "'56.55"
$ws.Range("E11").Value = "  +8.29%  "
$ws.Range("E12").Value = "  +17.36%  "
$ws.Range("D13").Value = "'0.0000275"
$ws.Range("E13").Value = "  +9.94%  "
$ws.Range("E14").Value = "  +7.41%  "
$ws.Range("D15").Value = "'4.104.41"
$ws.Range("E15").Value = "  +9.98%  "
$ws.Range("D16").Value = "'3.551.05"
$ws.Range("E16").Value = "  +10.43%  "
$ws.Range("D17").Value = "'67.857.26"
$ws.Range("E17").Value = "  +7.99%  "
$ws.Range("E18").Value = "  +5.98%  "
$ws.Range("E19").Value = "  +7.90%  "
$ws.Range("D20").Value = "'11.97"
$ws.Range("E20").Value = "  +9.86%  "
$ws.Range("E21").Value = "  +5.61%  "
$ws.Range("D22").Value = "'408.70"
$ws.Range("E22").Value = "  +12.37%  "
$ws.Range("D23").Value = "'3.99"
$ws.Range("E23").Value = "  +7.96%  "
$ws.Range("D24").Value = "'4.27"
$ws.Range("E24").Value = "  +9.98%  "
$ws.Range("D25").Value = "'85.38"
$ws.Range("E25").Value = "  +6.95%  "
$ws.Range("D26").Value = "'11.50"
$ws.Range("E26").Value = "  +4.94%  "
$ws.Range("D27").Value = "'2.99"
$ws.Range("E27").Value = "  +15.53%  "
$ws.Range("D28").Value = "'6.15"
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("E29").Value = "  +7.67%  "
$ws.Range("D30").Value = "'8.87"
$ws.Range("E30").Value = "  +9.35%  "
$ws.Range("D31").Value = "'695.90"
$ws.Range("E31").Value = "  +7.06%  "
$ws.Range("D32").Value = "'30.71"
$ws.Range("E32").Value = "  +9.29%  "
$ws.Range("D33").Value = "'6.88"
$ws.Range("E33").Value = "  +10.40%  "
$ws.Range("E34").Value = "  +7.20%  "
$ws.Range("D35").Value = "'0.113"
$ws.Range("E35").Value = "  +8.82%  "
$ws.Range("D36").Value = "'60.73"
$ws.Range("E36").Value = "  +5.32%  "
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").Value = "'0.0₃0843"
$ws.Range("E37").Value = "  +21.93%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "'39.40"
$ws.Range("E38").Value = "  +8.62%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  +7.62%  "
$ws.Range("D41").Value = "'0.140"
$ws.Range("E41").Value = "  +15.02%  "
$ws.Range("D42").Value = "'3.38"
$ws.Range("E42").Value = "  +18.43%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").Value = "'3.00"
$ws.Range("E44").Value = "  +16.48%  "
$ws.Range("D45").Value = "'3.064.73"
$ws.Range("E45").Value = "  +7.52%  "
$ws.Range("D46").Value = "'2.71"
$ws.Range("E46").Value = "  +7.70%  "
$ws.Range("E47").Value = "  +9.72%  "
$ws.Range("D48").Value = "'3.32"
$ws.Range("E48").Value = "  +14.19%  "
$ws.Range("D49").Value = "'9.13"
$ws.Range("E49").Value = "  +20.91%  "
$ws.Range("E50").Value = "  +2.32%  "
$ws.Range("E51").Value = "  +7.65%  "
